$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the price/volume columns keep their original "text" storage
# (the source data intentionally stores numeric-looking values such as
# "1.001" or "30.809.23" as plain text, not as numbers).
$ws.Range("D2:E51").NumberFormat = "@"

$rows = @(
    @(2, 'Bitcoin', 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc', '30.809.23', '  -0.56%  '),
    @(3, 'Ethereum', 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth', '1.936.04', '  -0.86%  '),
    @(4, 'TetherUSD', 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt', '1.001', '  +0.16%  '),
    @(5, 'BNB', 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb', '243.23', '  -0.89%  '),
    @(6, 'USDC', 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc', '0.9999', '  +0.07%  '),
    @(7, 'XRP', 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp', '0.4881', '  +0.04%  '),
    @(8, 'Cardano', 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada', '0.2943', '  -0.95%  '),
    @(9, 'Dogecoin', 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge', '0.06876', '  +0.35%  '),
    @(10, 'Solana', 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol', '19.27', '  +0.87%  '),
    @(11, 'Litecoin', 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc', '104.65', '  -2.66%  '),
    @(12, 'TRON', 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx', '0.07791', '  +0.74%  '),
    @(13, 'WrappedEther', 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth', '1.937.42', '  -0.97%  '),
    @(14, 'Polkadot', 'https://coinranking.com/coin/25W7FG7om+polkadot-dot', '5.334', '  -2.50%  '),
    @(15, 'Polygon', 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic', '0.6996', '  -1.15%  '),
    @(16, 'BitcoinCash', 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch', '272.84', '  -3.19%  '),
    @(17, 'WrappedBTC', 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc', '30.810.91', '  -0.74%  '),
    @(18, 'ShibaInu', 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib', '0.000007717', '  -0.53%  '),
    @(19, 'Avalanche', 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax', '13.08', '  -1.47%  '),
    @(20, 'Uniswap', 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni', '5.612', '  +1.47%  '),
    @(21, 'Dai', 'https://coinranking.com/coin/MoTuySvg7+dai-dai', '1.000', '  +0.10%  '),
    @(22, 'WrappedliquidstakedEther2.0', 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth', '2.192.84', '  -0.74%  '),
    @(23, 'BinanceUSD', 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd', '1.000', '  +0.21%  '),
    @(24, 'Chainlink', 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link', '6.523', '  +0.17%  '),
    @(25, 'Cosmos', 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom', '9.813', '  -0.23%  '),
    @(26, 'Monero', 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr', '165.20', '  -2.44%  '),
    @(27, 'EthereumClassic', 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc', '19.62', '  -1.91%  '),
    @(28, 'LidoDAOToken', 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo', '2.156', '  -3.38%  '),
    @(29, 'Stellar', 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm', '0.1037', '  -1.71%  '),
    @(30, 'Toncoin', 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton', '1.390', '  -2.30%  '),
    @(31, 'Filecoin', 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil', '4.584', '  +0.09%  '),
    @(32, 'PancakeSwap', 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake', '1.555', '  -2.02%  '),
    @(33, 'InternetComputer(DFINITY)', 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp', '4.378', '  -2.47%  '),
    @(34, 'Hedera', 'https://coinranking.com/coin/jad286TjB+hedera-hbar', '0.04883', '  -2.06%  '),
    @(35, 'ImmutableX', 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx', '0.7576', '  -0.35%  '),
    @(36, 'ARBITRUM', 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb', '1.148', '  -2.81%  '),
    @(37, 'Frax', 'https://coinranking.com/coin/KfWtaeV1W+frax-frax', '0.9994', '  +0.03%  '),
    @(38, 'HuobiToken', 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht', '2.730', '  +0.16%  '),
    @(39, 'VeChain', 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet', '0.02007', '  -1.20%  '),
    @(40, 'Aave', 'https://coinranking.com/coin/ixgUfzmLR+aave-aave', '80.06', '  +7.65%  '),
    @(41, 'MXToken', 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx', '2.656', '  -1.76%  '),
    @(42, 'FraxShare', 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs', '6.491', '  -0.22%  '),
    @(43, 'RenderToken', 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr', '2.082', '  -3.94%  '),
    @(44, 'TrustWalletToken', 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt', '0.9031', '  +2.08%  '),
    @(45, 'TheSandbox', 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand', '0.4428', '  -1.75%  '),
    @(46, 'Quant', 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt', '108.07', '  -1.27%  '),
    @(47, 'PaxDollar', 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp', '0.9998', '  +0.10%  '),
    @(48, 'Aptos', 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt', '7.770', '  -4.39%  '),
    @(49, 'Maker', 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr', '1.004.15', '  +2.25%  '),
    @(50, 'Algorand', 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo', '0.1245', '  -1.85%  '),
    @(51, 'Elrond', 'https://coinranking.com/coin/omwkOTglq+elrond-egld', '36.08', '  +0.78%  ')

)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}
